$d = $word.ActiveDocument

# The document currently ends with a paragraph that holds only the
# _GoBack bookmark, right before the sectPr.
$bookmarkParaIndex = $d.Paragraphs.Count
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)

# 1) Insert a brand-new plain paragraph right before it with the note text.
$bookmarkPara.Range.InsertParagraphBefore()
$notePara = $d.Paragraphs.Item($bookmarkParaIndex)
$notePara.Range.Text = "fix lỗi khi không có InitialModel"

# The bookmark paragraph shifted down by one.
$bookmarkParaIndex = $bookmarkParaIndex + 1
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)

# 2) Give the bookmark paragraph a leading run of spaces (before the
#    bookmark) and a trailing "Add-Migration InitialModel" run (after the
#    bookmark).
$bookmarkPara.Range.InsertAfter("Add-Migration InitialModel")
$bookmarkPara.Range.InsertBefore("              ")

# 3) Append a fresh empty paragraph right after the bookmark paragraph
#    before applying direct formatting, so the new paragraph does not pick
#    up the bold/red formatting.
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)
$bookmarkPara.Range.InsertParagraphAfter()

# 4) Now bold + color the bookmark paragraph's whole range (this also
#    stamps the paragraph-mark run properties, i.e. w:pPr/w:rPr).
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)
$bookmarkPara.Range.Font.Bold = $true
$bookmarkPara.Range.Font.Color = 255
